$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# --- Stable source cells used to copy the "N/A" placeholder text + style ---
# D14 = s14/text "0" ;  E14 = s14/text "***.*"  (both left untouched by this edit)

# --- Numeric cell updates ---
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("I14").Value = 1
$ws.Range("I14").NumberFormat = "#,##0"
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -66.666666666666
$ws.Range("M14").Value = -80
$ws.Range("N14").Value = -85.714285714285
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -60
$ws.Range("M15").Value = -60
$ws.Range("N15").Value = -90.47619047619
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -27.272727272727
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -35
$ws.Range("L16").Value = 4
$ws.Range("M16").Value = -53.571428571428
$ws.Range("N16").Value = -89.344262295082
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 41.176470588235
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 51
$ws.Range("K17").Value = 17.647058823529
$ws.Range("L17").Value = 3.448275862068
$ws.Range("M17").Value = -1.639344262295
$ws.Range("N17").Value = -60.264900662251
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("I18").Value = 20
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -52.380952380952
$ws.Range("L18").Value = -55.555555555555
$ws.Range("M18").Value = -44.444444444444
$ws.Range("N18").Value = -84.73282442748
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = -5
$ws.Range("I19").Value = 56
$ws.Range("J19").Value = 73
$ws.Range("K19").Value = -23.287671232876
$ws.Range("L19").Value = -25.333333333333
$ws.Range("M19").Value = -1.754385964912
$ws.Range("N19").Value = -20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -26.086956521739
$ws.Range("L20").Value = -55.263157894736
$ws.Range("M20").Value = -15
$ws.Range("N20").Value = -87.31343283582
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -20.253164556962
$ws.Range("I21").Value = 182
$ws.Range("J21").Value = 235
$ws.Range("K21").Value = -22.553191489361
$ws.Range("L21").Value = -27.490039840637
$ws.Range("M21").Value = -24.166666666666
$ws.Range("N21").Value = -75.98944591029
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("I22").Value = 3
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = -25
$ws.Range("M22").Value = 50
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = 17.647058823529
$ws.Range("L23").Value = 25
$ws.Range("M23").Value = 17.647058823529
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 16.666666666666
$ws.Range("I24").Value = 145
$ws.Range("J24").Value = 162
$ws.Range("K24").Value = -10.493827160493
$ws.Range("L24").Value = 10.687022900763
$ws.Range("M24").Value = 1.398601398601
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 400
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = -56.923076923076
$ws.Range("L25").Value = 40
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 166.666666666667
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 3.846153846153
$ws.Range("I26").Value = 65
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = -18.75
$ws.Range("L26").Value = -14.473684210526
$ws.Range("M26").Value = -53.900709219858
$ws.Range("D27").Value = 4
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -83.333333333333
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -44.444444444444
$ws.Range("C28").Value = 3
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 57.142857142857
$ws.Range("C29").Value = 3
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 3
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("I29").Value = 3
$ws.Range("I29").NumberFormat = "#,##0"
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -40
$ws.Range("M29").Value = -72.727272727272
$ws.Range("N29").Value = -92.682926829268
$ws.Range("C30").Value = 2
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 2
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 2
$ws.Range("I30").NumberFormat = "#,##0"
$ws.Range("K30").Value = -50
$ws.Range("L30").Value = -60
$ws.Range("M30").Value = -80
$ws.Range("N30").Value = -94.594594594594

# --- Text ("N/A"-placeholder) cell updates: copy value+style from a stable source cell ---
$ws.Range("D14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("D14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("D14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))
